# Split the long "Programa" (PT), "Programa" (EN/italic) and
# "Bibliografia" paragraphs into multiple <w:t> runs joined by
# manual line breaks (<w:br/>), matching the numbered-item structure
# of each paragraph.

$d = $word.ActiveDocument

function Insert-LineBreak {
    param(
        [string]$FindText,
        [string]$ReplaceText
    )
    $range = $d.Content
    $ok = $range.Find.Execute($FindText, $true, $false, $false, $false, $false, `
                               $true, 1, $false, $ReplaceText, 2)
    if (-not $ok) {
        throw "Find.Execute failed to locate: $FindText"
    }
}

# --- Portuguese "Programa" paragraph -----------------------------------
Insert-LineBreak "industriais.2. Processos" "industriais.^l2. Processos"
Insert-LineBreak "desidratados.4. Principais" "desidratados.^l4. Principais"
Insert-LineBreak "bioquímicas.5. Discussão" "bioquímicas.^l5. Discussão"
Insert-LineBreak "setores.6. Bioenergia" "setores.^l6. Bioenergia"

# --- English (italic) "Programa" paragraph ------------------------------
Insert-LineBreak "sectors.2. Biochemical" "sectors.^l2. Biochemical"
Insert-LineBreak "modifications3. Biochemical" "modifications^l3. Biochemical"
Insert-LineBreak "products.4. Main" "products.^l4. Main"
Insert-LineBreak "changes.5. Discussion" "changes.^l5. Discussion"
Insert-LineBreak "sectors.6. Bioenergy" "sectors.^l6. Bioenergy"

# --- "Bibliografia" paragraph --------------------------------------------
Insert-LineBreak "1.GAVA" "1.^lGAVA"
Insert-LineBreak "9788521313823.2.LIMA" "9788521313823.^l2.^lLIMA"
Insert-LineBreak "9788521214571.3.Moraes" "9788521214571.^l3.^lMoraes"
